{"js": "// Update the worksheet date and all 25 \"two-digit \u00d7 two-digit\" answer\n// cells. Each old value is unique in the document at the time it is\n// searched for, so a sequence of single-match search/replace calls\n// (applied in the same top-to-bottom, left-to-right order as the XML\n// diff) reproduces the edit deterministically \u2014 including the one spot\n// (55\u00d774=4070 -> 14\u00d735=490) where a freshly-written value happens to\n// match an *original* value that was already rewritten earlier in the\n// pass.\nconst replacements = [\n  [\"2024-08-12 Monday\", \"2024-08-13 Tuesday\"],\n  [\"82\u00d799=8118\", \"50\u00d712=600\"],\n  [\"96\u00d731=2976\", \"27\u00d742=1134\"],\n  [\"27\u00d773=1971\", \"27\u00d732=864\"],\n  [\"99\u00d794=9306\", \"28\u00d721=588\"],\n  [\"43\u00d720=860\", \"63\u00d738=2394\"],\n  [\"14\u00d735=490\", \"60\u00d745=2700\"],\n  [\"60\u00d753=3180\", \"49\u00d779=3871\"],\n  [\"50\u00d758=2900\", \"34\u00d740=1360\"],\n  [\"20\u00d778=1560\", \"30\u00d749=1470\"],\n  [\"45\u00d778=3510\", \"74\u00d747=3478\"],\n  [\"53\u00d720=1060\", \"63\u00d748=3024\"],\n  [\"55\u00d774=4070\", \"14\u00d735=490\"],\n  [\"28\u00d756=1568\", \"87\u00d731=2697\"],\n  [\"32\u00d749=1568\", \"29\u00d741=1189\"],\n  [\"50\u00d760=3000\", \"25\u00d754=1350\"],\n  [\"46\u00d765=2990\", \"67\u00d732=2144\"],\n  [\"53\u00d748=2544\", \"74\u00d763=4662\"],\n  [\"56\u00d757=3192\", \"71\u00d738=2698\"],\n  [\"23\u00d733=759\", \"95\u00d746=4370\"],\n  [\"25\u00d790=2250\", \"17\u00d797=1649\"],\n  [\"37\u00d783=3071\", \"22\u00d776=1672\"],\n  [\"74\u00d743=3182\", \"35\u00d750=1750\"],\n  [\"32\u00d764=2048\", \"48\u00d778=3744\"],\n  [\"50\u00d740=2000\", \"40\u00d741=1640\"],\n  [\"27\u00d768=1836\", \"91\u00d780=7280\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 \"two-digit x two-digit\" answer\n# cells using Find/Replace on the document's Content range. Each old\n# value is unique in the document at the time it is searched for, so\n# running the replacements in the same top-to-bottom, left-to-right\n# order as the XML diff reproduces the edit deterministically --\n# including the one spot (55x74=4070 -> 14x35=490) where a freshly\n# written value happens to match an *original* value that was already\n# rewritten earlier in the pass.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-12 Monday\", \"2024-08-13 Tuesday\"),\n    @(\"82\u00d799=8118\", \"50\u00d712=600\"),\n    @(\"96\u00d731=2976\", \"27\u00d742=1134\"),\n    @(\"27\u00d773=1971\", \"27\u00d732=864\"),\n    @(\"99\u00d794=9306\", \"28\u00d721=588\"),\n    @(\"43\u00d720=860\", \"63\u00d738=2394\"),\n    @(\"14\u00d735=490\", \"60\u00d745=2700\"),\n    @(\"60\u00d753=3180\", \"49\u00d779=3871\"),\n    @(\"50\u00d758=2900\", \"34\u00d740=1360\"),\n    @(\"20\u00d778=1560\", \"30\u00d749=1470\"),\n    @(\"45\u00d778=3510\", \"74\u00d747=3478\"),\n    @(\"53\u00d720=1060\", \"63\u00d748=3024\"),\n    @(\"55\u00d774=4070\", \"14\u00d735=490\"),\n    @(\"28\u00d756=1568\", \"87\u00d731=2697\"),\n    @(\"32\u00d749=1568\", \"29\u00d741=1189\"),\n    @(\"50\u00d760=3000\", \"25\u00d754=1350\"),\n    @(\"46\u00d765=2990\", \"67\u00d732=2144\"),\n    @(\"53\u00d748=2544\", \"74\u00d763=4662\"),\n    @(\"56\u00d757=3192\", \"71\u00d738=2698\"),\n    @(\"23\u00d733=759\", \"95\u00d746=4370\"),\n    @(\"25\u00d790=2250\", \"17\u00d797=1649\"),\n    @(\"37\u00d783=3071\", \"22\u00d776=1672\"),\n    @(\"74\u00d743=3182\", \"35\u00d750=1750\"),\n    @(\"32\u00d764=2048\", \"48\u00d778=3744\"),\n    @(\"50\u00d740=2000\", \"40\u00d741=1640\"),\n    @(\"27\u00d768=1836\", \"91\u00d780=7280\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
